$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '61.859.66', '  -0.59%  ')
    3 = @('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.418.77', '  -0.62%  ')
    4 = @('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  +0.12%  ')
    5 = @('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '410.02', '  +0.05%  ')
    6 = @('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '129.29', '  -0.57%  ')
    7 = @('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.631', '  -0.52%  ')
    8 = @('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  +0.00%  ')
    9 = @('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.733', '  -3.79%  ')
    10 = @('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.139', '  -1.66%  ')
    11 = @('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '43.30', '  +0.30%  ')
    12 = @('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000222', '  +16.42%  ')
    13 = @('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '9.27', '  +5.29%  ')
    14 = @('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.959.91', '  -0.58%  ')
    15 = @('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.141', '  +0.32%  ')
    16 = @('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '21.20', '  +3.67%  ')
    17 = @('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.403.67', '  -1.03%  ')
    18 = @('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '12.32', '  +8.21%  ')
    19 = @('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.08', '  +2.74%  ')
    20 = @('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '61.784.68', '  -0.66%  ')
    21 = @('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '501.56', '  +29.15%  ')
    22 = @('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '92.00', '  +3.97%  ')
    23 = @('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '3.32', '  +4.32%  ')
    24 = @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '13.45', '  -0.11%  ')
    25 = @('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.34', '  +3.69%  ')
    26 = @('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '34.64', '  +8.40%  ')
    27 = @('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '9.29', '  +9.05%  ')
    28 = @('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '7.59', '  -1.26%  ')
    29 = @('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '12.17', '  +2.91%  ')
    30 = @('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.69', '  -1.18%  ')
    31 = @('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.114', '  -1.71%  ')
    32 = @('Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.168', '  -2.30%  ')
    33 = @('InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '41.91', '  -5.11%  ')
    34 = @('OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '59.53', '  +13.91%  ')
    35 = @('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  -0.01%  ')
    36 = @('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0499', '  +1.07%  ')
    37 = @('FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '0.998', '  +0.06%  ')
    38 = @('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '3.46', '  +2.74%  ')
    39 = @('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.137', '  +3.54%  ')
    40 = @('WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '2.72', '  +17.02%  ')
    41 = @('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '146.89', '  +3.42%  ')
    42 = @('TheGraph', 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt', '0.319', '  +2.11%  ')
    43 = @('Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '2.92', '  -0.09%  ')
    44 = @('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '2.10', '  +6.35%  ')
    45 = @('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '4.34', '  +8.40%  ')
    46 = @('Celestia', 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia', '16.69', '  -0.37%  ')
    47 = @('ThetaToken', 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta', '2.31', '  +19.83%  ')
    48 = @('BitcoinSV', 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv', '119.45', '  +27.11%  ')
    49 = @('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '22.84', '  +4.48%  ')
    50 = @('Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.144', '  +17.49%  ')
    51 = @('Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '2.139.81', '  +0.87%  ')
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]

    # Column D must stay plain text even when it looks numeric,
    # so force Text format, assign, then restore the original style
    # to avoid leaving a stray number-format style behind.
    $dCell = $ws.Cells.Item($r, 4)
    $origStyle = $dCell.Style
    $dCell.NumberFormat = "@"
    $dCell.Value = $vals[2]
    $dCell.Style = $origStyle

    $ws.Cells.Item($r, 5).Value = $vals[3]
}

